$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price text would otherwise be auto-parsed by Excel as a
# number (single decimal point, e.g. "505.04"). For these we briefly force
# the built-in Text format so the literal string is stored, then clear the
# formatting again so no visible/style change is left behind (matches the
# other price cells, which carry no explicit cell style).
$textFormatCells = @("D5", "D6", "D8", "D10", "D12", "D16", "D19", "D21", "D22", "D24", "D25", "D26", "D28", "D29", "D30", "D33", "D34", "D37", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D51")

# New values scraped for this run, keyed by cell address.
$newValues = [ordered]@{
    "D2" = "56.543.28"
    "E2" = "  -0.07%  "
    "D3" = "2.385.60"
    "E3" = "  +0.23%  "
    "E4" = "  +0.01%  "
    "D5" = "505.04"
    "E5" = "  +1.14%  "
    "D6" = "133.21"
    "E6" = "  +3.43%  "
    "E7" = "  +0.16%  "
    "D8" = "0.551"
    "E8" = "  +0.20%  "
    "D9" = "2.386.39"
    "E9" = "  -0.79%  "
    "D10" = "0.0975"
    "E10" = "  +2.19%  "
    "E11" = "  +0.51%  "
    "D12" = "0.332"
    "E12" = "  +4.22%  "
    "E13" = "  +0.22%  "
    "D14" = "2.808.28"
    "E14" = "  +0.10%  "
    "D15" = "56.490.95"
    "E15" = "  +0.08%  "
    "D16" = "21.69"
    "E16" = "  +0.81%  "
    "E17" = "  +0.85%  "
    "D18" = "2.389.18"
    "E18" = "  -1.72%  "
    "D19" = "10.17"
    "E19" = "  +0.34%  "
    "E20" = "  +0.28%  "
    "D21" = "309.23"
    "E21" = "  -0.16%  "
    "D22" = "6.25"
    "E22" = "  +0.26%  "
    "E23" = "  +0.35%  "
    "D24" = "5.61"
    "E24" = "  -4.31%  "
    "D25" = "65.41"
    "E25" = "  +0.31%  "
    "D26" = "0.998"
    "E26" = "  +0.07%  "
    "E27" = "  -0.17%  "
    "D28" = "0.369"
    "E28" = "  -1.84%  "
    "D29" = "7.32"
    "E29" = "  +1.27%  "
    "D30" = "174.28"
    "E30" = "  +0.78%  "
    "D31" = "0.0₃0725"
    "E31" = "  +1.62%  "
    "E32" = "  -0.33%  "
    "D33" = "1.11"
    "E33" = "  +1.49%  "
    "D34" = "5.85"
    "E34" = "  -4.24%  "
    "E35" = "  +0.13%  "
    "E36" = "  +0.17%  "
    "D37" = "17.76"
    "E37" = "  -0.10%  "
    "E38" = "  -0.50%  "
    "E39" = "  +1.41%  "
    "D40" = "36.70"
    "E40" = "  +2.37%  "
    "D41" = "0.812"
    "E41" = "  +2.65%  "
    "D42" = "1.43"
    "E42" = "  +0.26%  "
    "D43" = "131.72"
    "E43" = "  +0.66%  "
    "D44" = "3.38"
    "E44" = "  +1.06%  "
    "D45" = "4.82"
    "E45" = "  +0.99%  "
    "D46" = "0.566"
    "E46" = "  -1.49%  "
    "D47" = "0.0908"
    "E47" = "  +1.17%  "
    "D48" = "246.06"
    "E48" = "  -2.81%  "
    "D49" = "0.0484"
    "E49" = "  -0.04%  "
    "E50" = "  +1.16%  "
    "D51" = "17.10"
    "E51" = "  +6.42%  "
}

foreach ($addr in $newValues.Keys) {
    $cell = $ws.Range($addr)
    if ($textFormatCells -contains $addr) {
        $cell.NumberFormat = "@"
        $cell.Value = $newValues[$addr]
        $cell.ClearFormats()
    } else {
        $cell.Value = $newValues[$addr]
    }
}
